$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '22.323.33'
Set-TextCell $ws 'E2' '  -0.71%  '
Set-TextCell $ws 'D3' '1.565.30'
Set-TextCell $ws 'E3' '  -0.50%  '
Set-TextCell $ws 'E4' '  +0.09%  '
Set-TextCell $ws 'D5' '1.002'
Set-TextCell $ws 'E5' '  +0.07%  '
Set-TextCell $ws 'D6' '286.57'
Set-TextCell $ws 'D7' '0.3757'
Set-TextCell $ws 'E7' '  +2.62%  '
Set-TextCell $ws 'D8' '0.3278'
Set-TextCell $ws 'E8' '  -1.83%  '
Set-TextCell $ws 'D9' '45.49'
Set-TextCell $ws 'E9' '  -5.44%  '
Set-TextCell $ws 'D10' '1.142'
Set-TextCell $ws 'E10' '  +1.09%  '
Set-TextCell $ws 'D11' '0.07399'
Set-TextCell $ws 'E11' '  -0.60%  '
Set-TextCell $ws 'D12' '1.003'
Set-TextCell $ws 'E12' '  +0.12%  '
Set-TextCell $ws 'D13' '20.36'
Set-TextCell $ws 'E13' '  -2.57%  '
Set-TextCell $ws 'D14' '5.840'
Set-TextCell $ws 'E14' '  -2.61%  '
Set-TextCell $ws 'D15' '6.801'
Set-TextCell $ws 'E15' '  -1.71%  '
Set-TextCell $ws 'D16' '1.561.62'
Set-TextCell $ws 'E16' '  -0.78%  '
Set-TextCell $ws 'D17' '0.00001093'
Set-TextCell $ws 'E17' '  -1.63%  '
Set-TextCell $ws 'D18' '0.06711'
Set-TextCell $ws 'E18' '  -0.60%  '
Set-TextCell $ws 'D19' '86.04'
Set-TextCell $ws 'E19' '  -2.29%  '
Set-TextCell $ws 'D20' '1.002'
Set-TextCell $ws 'E20' '  +0.01%  '
Set-TextCell $ws 'D21' '6.347'
Set-TextCell $ws 'E21' '  -0.69%  '
Set-TextCell $ws 'D22' '16.25'
Set-TextCell $ws 'E22' '  -1.24%  '
Set-TextCell $ws 'E23' '  -3.64%  '
Set-TextCell $ws 'D24' '22.322.91'
Set-TextCell $ws 'E24' '  -0.69%  '
Set-TextCell $ws 'E25' '  -3.82%  '
Set-TextCell $ws 'D26' '2.519'
Set-TextCell $ws 'E26' '  -3.84%  '
Set-TextCell $ws 'D27' '149.99'
Set-TextCell $ws 'E27' '  -1.66%  '
Set-TextCell $ws 'D28' '19.42'
Set-TextCell $ws 'E28' '  -1.12%  '
Set-TextCell $ws 'D29' '4.885'
Set-TextCell $ws 'E29' '  -2.59%  '
Set-TextCell $ws 'D30' '123.63'
Set-TextCell $ws 'E30' '  -0.47%  '
Set-TextCell $ws 'D31' '1.740.73'
Set-TextCell $ws 'E31' '  -0.60%  '
Set-TextCell $ws 'D32' '1.049'
Set-TextCell $ws 'E32' '  +0.91%  '
Set-TextCell $ws 'D33' '5.911'
Set-TextCell $ws 'E33' '  -4.31%  '
Set-TextCell $ws 'D34' '1.906'
Set-TextCell $ws 'E34' '  -4.27%  '
Set-TextCell $ws 'D35' '9.470'
Set-TextCell $ws 'E35' '  -3.14%  '
Set-TextCell $ws 'D36' '0.08222'
Set-TextCell $ws 'E36' '  -0.62%  '
Set-TextCell $ws 'D37' '0.02379'
Set-TextCell $ws 'E37' '  -2.34%  '
Set-TextCell $ws 'B38' 'Hedera'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D38' '0.06287'
Set-TextCell $ws 'E38' '  -3.14%  '
Set-TextCell $ws 'B39' 'TrustWalletToken'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D39' '1.278'
Set-TextCell $ws 'E39' '  -1.83%  '
Set-TextCell $ws 'E40' '  -3.84%  '
Set-TextCell $ws 'D41' '5.248'
Set-TextCell $ws 'E41' '  -3.24%  '
Set-TextCell $ws 'D42' '11.00'
Set-TextCell $ws 'E42' '  -2.92%  '
Set-TextCell $ws 'D43' '0.6070'
Set-TextCell $ws 'E43' '  -3.47%  '
Set-TextCell $ws 'D44' '1.002'
Set-TextCell $ws 'E44' '  +0.06%  '
Set-TextCell $ws 'D45' '13.69'
Set-TextCell $ws 'E45' '  -2.19%  '
Set-TextCell $ws 'D46' '3.741'
Set-TextCell $ws 'E46' '  -0.36%  '
Set-TextCell $ws 'D47' '0.5886'
Set-TextCell $ws 'E47' '  -3.00%  '
Set-TextCell $ws 'B48' 'Quant'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws 'D48' '123.91'
Set-TextCell $ws 'E48' '  -0.63%  '
Set-TextCell $ws 'B49' 'NEARProtocol'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D49' '1.992'
Set-TextCell $ws 'E49' '  -2.91%  '
Set-TextCell $ws 'D50' '1.178'
Set-TextCell $ws 'E50' '  -3.63%  '
Set-TextCell $ws 'D51' '0.07138'
Set-TextCell $ws 'E51' '  -1.15%  '
